# Atualizado por script em 07-11-2023 08:45
#
# Two kinds of changes:
#  1. Four pairs of adjacent match rows had their match-detail columns
#     (F:V - teams, scores, odds, timestamps, url) swapped between the two
#     rows in the pair (the "Indice"/date columns A:E stay put).
#  2. Two brand-new match rows (232, 233) are appended at the end of the
#     sheet, pushing the dimension out to A1:V233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap F:V content between row 53 and row 54
# ---------------------------------------------------------------------
$ws.Range("F53").Value = "Brommapojkarna"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Goteborg"
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 2.74
$ws.Range("K53").Value = "04/05/2023 19:12"
$ws.Range("L53").Value = 2.87
$ws.Range("M53").Value = "08/05/2023 18:59"
$ws.Range("N53").Value = 3.33
$ws.Range("O53").Value = "04/05/2023 19:12"
$ws.Range("P53").Value = 3.66
$ws.Range("Q53").Value = "08/05/2023 18:59"
$ws.Range("R53").Value = 2.68
$ws.Range("S53").Value = "04/05/2023 19:12"
$ws.Range("T53").Value = 2.47
$ws.Range("U53").Value = "08/05/2023 18:56"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/brommapojkarna-ifk-goteborg/lbbxFdBs/"

$ws.Range("F54").Value = "Degerfors"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = "Norrkoping"
$ws.Range("I54").Value = 2
$ws.Range("J54").Value = 2.77
$ws.Range("K54").Value = "04/05/2023 19:12"
$ws.Range("L54").Value = 2.87
$ws.Range("M54").Value = "08/05/2023 18:54"
$ws.Range("N54").Value = 3.47
$ws.Range("O54").Value = "04/05/2023 19:12"
$ws.Range("P54").Value = 3.47
$ws.Range("Q54").Value = "08/05/2023 18:59"
$ws.Range("R54").Value = 2.62
$ws.Range("S54").Value = "04/05/2023 19:12"
$ws.Range("T54").Value = 2.56
$ws.Range("U54").Value = "08/05/2023 18:59"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/degerfors-norrkoping/rLctEGQm/"

# ---------------------------------------------------------------------
# 2) Swap F:V content between row 196 and row 197
# ---------------------------------------------------------------------
$ws.Range("F196").Value = "Hacken"
$ws.Range("G196").Value = 2
$ws.Range("H196").Value = "AIK"
$ws.Range("I196").Value = 0
$ws.Range("J196").Value = 1.74
$ws.Range("K196").Value = "24/09/2023 16:42"
$ws.Range("L196").Value = 1.92
$ws.Range("M196").Value = "01/10/2023 14:55"
$ws.Range("N196").Value = 4.2
$ws.Range("O196").Value = "24/09/2023 16:42"
$ws.Range("P196").Value = 4.07
$ws.Range("Q196").Value = "01/10/2023 14:58"
$ws.Range("R196").Value = 4.51
$ws.Range("S196").Value = "24/09/2023 16:42"
$ws.Range("T196").Value = 3.82
$ws.Range("U196").Value = "01/10/2023 14:48"
$ws.Range("V196").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/hacken-aik/n5Cb3ZVB/"

$ws.Range("F197").Value = "Sirius"
$ws.Range("G197").Value = 3
$ws.Range("H197").Value = "Degerfors"
$ws.Range("I197").Value = 2
$ws.Range("J197").Value = 1.62
$ws.Range("K197").Value = "25/09/2023 18:12"
$ws.Range("L197").Value = 1.62
$ws.Range("M197").Value = "01/10/2023 14:48"
$ws.Range("N197").Value = 4.4
$ws.Range("O197").Value = "25/09/2023 18:12"
$ws.Range("P197").Value = 4.31
$ws.Range("Q197").Value = "01/10/2023 14:52"
$ws.Range("R197").Value = 5.05
$ws.Range("S197").Value = "25/09/2023 18:12"
$ws.Range("T197").Value = 5.5
$ws.Range("U197").Value = "01/10/2023 14:56"
$ws.Range("V197").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/sirius-degerfors/QeB22goI/"

# ---------------------------------------------------------------------
# 3) Swap F:V content between row 220 and row 221
# ---------------------------------------------------------------------
$ws.Range("F220").Value = "Mjallby"
$ws.Range("G220").Value = 2
$ws.Range("H220").Value = "Hacken"
$ws.Range("I220").Value = 1
$ws.Range("J220").Value = 3.31
$ws.Range("K220").Value = "22/10/2023 20:15"
$ws.Range("L220").Value = 4.35
$ws.Range("M220").Value = "29/10/2023 14:51"
$ws.Range("N220").Value = 3.68
$ws.Range("O220").Value = "22/10/2023 20:15"
$ws.Range("P220").Value = 4.21
$ws.Range("Q220").Value = "29/10/2023 14:51"
$ws.Range("R220").Value = 2.18
$ws.Range("S220").Value = "22/10/2023 20:15"
$ws.Range("T220").Value = 1.78
$ws.Range("U220").Value = "29/10/2023 14:51"
$ws.Range("V220").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/mjallby-hacken/IVVkNBxi/"

$ws.Range("F221").Value = "Norrkoping"
$ws.Range("G221").Value = 0
$ws.Range("H221").Value = "Malmo FF"
$ws.Range("I221").Value = 1
$ws.Range("J221").Value = 4.05
$ws.Range("K221").Value = "22/10/2023 20:15"
$ws.Range("L221").Value = 7.1
$ws.Range("M221").Value = "29/10/2023 14:56"
$ws.Range("N221").Value = 3.9
$ws.Range("O221").Value = "22/10/2023 20:15"
$ws.Range("P221").Value = 5
$ws.Range("Q221").Value = "29/10/2023 14:56"
$ws.Range("R221").Value = 1.88
$ws.Range("S221").Value = "22/10/2023 20:15"
$ws.Range("T221").Value = 1.45
$ws.Range("U221").Value = "29/10/2023 14:56"
$ws.Range("V221").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/norrkoping-malmo-ff/6ijgNTMp/"

# ---------------------------------------------------------------------
# 4) Swap F:V content between row 228 and row 229
# ---------------------------------------------------------------------
$ws.Range("F228").Value = "Halmstad"
$ws.Range("G228").Value = 3
$ws.Range("H228").Value = "Kalmar"
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = 2.38
$ws.Range("K228").Value = "29/10/2023 17:43"
$ws.Range("L228").Value = 2.53
$ws.Range("M228").Value = "05/11/2023 14:57"
$ws.Range("N228").Value = 3.32
$ws.Range("O228").Value = "29/10/2023 17:43"
$ws.Range("P228").Value = 3.33
$ws.Range("Q228").Value = "05/11/2023 14:47"
$ws.Range("R228").Value = 3.21
$ws.Range("S228").Value = "29/10/2023 17:43"
$ws.Range("T228").Value = 3.01
$ws.Range("U228").Value = "05/11/2023 14:57"
$ws.Range("V228").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/halmstad-kalmar/Cdwmaaj1/"

$ws.Range("F229").Value = "Hacken"
$ws.Range("G229").Value = 4
$ws.Range("H229").Value = "Malmo FF"
$ws.Range("I229").Value = 2
$ws.Range("J229").Value = 2.32
$ws.Range("K229").Value = "29/10/2023 15:13"
$ws.Range("L229").Value = 3.3
$ws.Range("M229").Value = "05/11/2023 14:59"
$ws.Range("N229").Value = 3.92
$ws.Range("O229").Value = "29/10/2023 15:13"
$ws.Range("P229").Value = 4.17
$ws.Range("Q229").Value = "05/11/2023 14:55"
$ws.Range("R229").Value = 2.89
$ws.Range("S229").Value = "29/10/2023 15:13"
$ws.Range("T229").Value = 2.06
$ws.Range("U229").Value = "05/11/2023 14:58"
$ws.Range("V229").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/hacken-malmo-ff/UXmr0uye/"

# ---------------------------------------------------------------------
# 5) Append two new match rows (232, 233) after the existing last row
#    (231). Copy number formats/styles from the row above first (A uses
#    the bold/bordered "index" style, E uses the datetime style), then
#    set the actual values so the copy doesn't clobber them.
#
#    D ("temporada") is the literal text "2023" (not a number) in every
#    other row, same as the row above it - a plain `.Value = "2023"`
#    assignment gets auto-coerced to a numeric cell, so instead it's
#    copied verbatim (value + format together) from D231, cell-by-cell
#    (pasting a single cell into a multi-cell destination range only
#    fills the first cell here, so each target cell is copied to on its
#    own).
# ---------------------------------------------------------------------
$ws.Range("A231").Copy()
$ws.Range("A232:A233").PasteSpecial(-4122)
$ws.Range("E231").Copy()
$ws.Range("E232:E233").PasteSpecial(-4122)
$ws.Range("D231").Copy()
$ws.Range("D232").PasteSpecial(-4104)
$ws.Range("D231").Copy()
$ws.Range("D233").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("A232").Value = 231
$ws.Range("B232").Value = "sweden"
$ws.Range("C232").Value = "allsvenskan"
$ws.Range("E232").Value = 45236.79166666666
$ws.Range("F232").Value = "Norrkoping"
$ws.Range("G232").Value = 4
$ws.Range("H232").Value = "Varberg"
$ws.Range("I232").Value = 3
$ws.Range("J232").Value = 1.49
$ws.Range("K232").Value = "30/10/2023 19:13"
$ws.Range("L232").Value = 1.44
$ws.Range("M232").Value = "06/11/2023 18:22"
$ws.Range("N232").Value = 4.75
$ws.Range("O232").Value = "30/10/2023 19:13"
$ws.Range("P232").Value = 5.21
$ws.Range("Q232").Value = "06/11/2023 18:22"
$ws.Range("R232").Value = 6.48
$ws.Range("S232").Value = "30/10/2023 19:13"
$ws.Range("T232").Value = 6.94
$ws.Range("U232").Value = "06/11/2023 18:22"
$ws.Range("V232").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/norrkoping-varberg/M7pz216r/"

$ws.Range("A233").Value = 232
$ws.Range("B233").Value = "sweden"
$ws.Range("C233").Value = "allsvenskan"
$ws.Range("E233").Value = 45236.79861111111
$ws.Range("F233").Value = "Goteborg"
$ws.Range("G233").Value = 1
$ws.Range("H233").Value = "AIK"
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = 2.55
$ws.Range("K233").Value = "30/10/2023 19:13"
$ws.Range("L233").Value = 2.4
$ws.Range("M233").Value = "06/11/2023 19:07"
$ws.Range("N233").Value = 3.32
$ws.Range("O233").Value = "30/10/2023 19:13"
$ws.Range("P233").Value = 3.13
$ws.Range("Q233").Value = "06/11/2023 19:07"
$ws.Range("R233").Value = 2.97
$ws.Range("S233").Value = "30/10/2023 19:13"
$ws.Range("T233").Value = 3.44
$ws.Range("U233").Value = "06/11/2023 19:07"
$ws.Range("V233").Value = "https://www.betexplorer.com/football/sweden/allsvenskan/ifk-goteborg-aik/prN6JTyG/"
